$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J, matching the style of the existing header row
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I0 and IF data for rows 2-63
$data = @{
    2 = @(8, 8)
    3 = @(7, 8)
    4 = @(8, 9)
    5 = @(7, 8)
    6 = @(8, 8)
    7 = @(9, 9)
    8 = @(10, 10)
    9 = @(8, 8)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(9, 9)
    15 = @(7, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(8, 8)
    19 = @(9, 9)
    20 = @(8, 8)
    21 = @(10, 10)
    22 = @(9, 9)
    23 = @(8, 8)
    24 = @(9, 9)
    25 = @(8, 8)
    26 = @(7, 8)
    27 = @(6, 6)
    28 = @(9, 9)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(7, 7)
    32 = @(8, 8)
    33 = @(8, 8)
    34 = @(8, 8)
    35 = @(7, 8)
    36 = @(6, 7)
    37 = @(9, 9)
    38 = @(7, 8)
    39 = @(8, 8)
    40 = @(8, 8)
    41 = @(9, 10)
    42 = @(8, 8)
    43 = @(6, 7)
    44 = @(9, 9)
    45 = @(8, 8)
    46 = @(8, 9)
    47 = @(10, 10)
    48 = @(9, 9)
    49 = @(9, 9)
    50 = @(9, 9)
    51 = @(8, 8)
    52 = @(9, 9)
    53 = @(9, 9)
    54 = @(9, 9)
    55 = @(9, 9)
    56 = @(9, 9)
    57 = @(9, 9)
    58 = @(8, 9)
    59 = @(9, 9)
    60 = @(9, 9)
    61 = @(7, 7)
    62 = @(9, 9)
    63 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
